# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values. Update them to reflect the regenerated
# (Strike -> K) calculation results for each saved round/row.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 13
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 4
$ws.Range("G11").Value = 1
